# "Update gh-pages to output generated at 456a3b4"
#
# The scraper re-ran and produced slightly different numbers for this
# commit:
#   - every start-date string in column B changes from dotted
#     (2024.MM.DD) to dashed (2024-MM-DD) notation, on every sheet that
#     has data rows (展览, 演出, 全部类型)
#   - a handful of "want to go" counters in column F ticked up by one
#   - the "已停售" (sold out) ticket-status label for the 绘时国乙1.0
#     event was corrected to "不可售" (not for sale), on 展览 and 全部类型

$wb = $excel.ActiveWorkbook

function Update-DateColumn($SheetName, $LastRow) {
    $ws = $wb.Worksheets.Item($SheetName)

    for ($r = 2; $r -le $LastRow; $r++) {
        $cell = $ws.Cells.Item($r, 2)  # column B = start date
        $text = $cell.Value()
        if ($text -ne $null -and $text -like "2024.*") {
            # Leading apostrophe forces plain text entry, the same way
            # typing it in the Excel UI would, so the dashed date is not
            # silently reinterpreted as a date serial number.
            $cell.Value = "'" + $text.Replace(".", "-")
        }
    }
}

Update-DateColumn "展览" 20
Update-DateColumn "演出" 2
Update-DateColumn "全部类型" 21

# 展览 sheet: bump interest counters, fix ticket-status label
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 53
$ws1.Cells.Item(5, 6).Value = 344
$ws1.Cells.Item(7, 6).Value = 10894
$ws1.Cells.Item(14, 6).Value = 12649
$ws1.Cells.Item(10, 7).Value = "不可售"

# 全部类型 sheet: same underlying rows, shifted down by one (it also
# includes the 演出 row at the top)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 53
$ws4.Cells.Item(6, 6).Value = 344
$ws4.Cells.Item(8, 6).Value = 10894
$ws4.Cells.Item(15, 6).Value = 12649
$ws4.Cells.Item(11, 7).Value = "不可售"
